$d = $word.ActiveDocument

$pairs = @(
    @{old="790×7="; new="865×6="},
    @{old="910×2="; new="253×9="},
    @{old="192×9="; new="827×9="},
    @{old="250×3="; new="391×5="},
    @{old="754×7="; new="742×7="},
    @{old="216×8="; new="533×3="},
    @{old="332×6="; new="804×8="},
    @{old="404×5="; new="946×7="},
    @{old="248×6="; new="357×2="},
    @{old="678×4="; new="931×5="},
    @{old="851×3="; new="328×8="},
    @{old="601×5="; new="398×4="},
    @{old="713×5="; new="845×8="},
    @{old="382×4="; new="730×2="},
    @{old="956×9="; new="251×8="},
    @{old="196×6="; new="890×8="},
    @{old="702×2="; new="357×3="},
    @{old="887×8="; new="204×9="},
    @{old="524×2="; new="746×2="},
    @{old="453×8="; new="590×5="},
    @{old="316×3="; new="813×3="},
    @{old="616×6="; new="634×6="},
    @{old="652×3="; new="508×7="},
    @{old="606×9="; new="626×6="},
    @{old="654×6="; new="469×4="}
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $pair.new, 2)
}
